# Applies the cryptos-list price/volume/ranking refresh described by the commit
# "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Test-ExactText($actual, $expected) {
    if ($actual.Length -ne $expected.Length) { return $false }
    for ($i = 0; $i -lt $expected.Length; $i++) {
        if ($actual[$i] -ne $expected[$i]) { return $false }
    }
    return $true
}

function Set-ExactText($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.Value = $value
    if (-not (Test-ExactText $cell.Text $value)) {
        $cell.Value = "'" + $value
        $cell.Style = "Normal"
    }
}

Set-ExactText "D2" "63.071.81"
Set-ExactText "E2" "  -1.85%  "
Set-ExactText "D3" "3.125.51"
Set-ExactText "E3" "  -0.46%  "
Set-ExactText "E4" "  -0.05%  "
Set-ExactText "D5" "595.41"
Set-ExactText "E5" "  -2.44%  "
Set-ExactText "D6" "136.47"
Set-ExactText "E6" "  -4.90%  "
Set-ExactText "E7" "  -0.07%  "
Set-ExactText "D8" "3.118.59"
Set-ExactText "E8" "  -0.54%  "
Set-ExactText "E9" "  -2.45%  "
Set-ExactText "D10" "0.146"
Set-ExactText "E10" "  -3.45%  "
Set-ExactText "D11" "5.20"
Set-ExactText "E11" "  -3.96%  "
Set-ExactText "D12" "0.460"
Set-ExactText "E12" "  -3.62%  "
Set-ExactText "D13" "0.0000248"
Set-ExactText "E13" "  -2.90%  "
Set-ExactText "D14" "34.26"
Set-ExactText "E14" "  -3.61%  "
Set-ExactText "D15" "3.633.97"
Set-ExactText "E15" "  -0.61%  "
Set-ExactText "E16" "  +1.48%  "
Set-ExactText "D17" "63.008.57"
Set-ExactText "E17" "  -2.01%  "
Set-ExactText "D18" "3.119.12"
Set-ExactText "E18" "  -0.88%  "
Set-ExactText "D19" "6.74"
Set-ExactText "E19" "  -2.07%  "
Set-ExactText "D20" "477.12"
Set-ExactText "E20" "  -0.12%  "
Set-ExactText "D21" "14.18"
Set-ExactText "E21" "  -3.72%  "
Set-ExactText "D22" "0.698"
Set-ExactText "E22" "  -3.62%  "
Set-ExactText "D23" "7.69"
Set-ExactText "E23" "  -1.86%  "
Set-ExactText "D24" "87.48"
Set-ExactText "E24" "  +2.43%  "
Set-ExactText "D25" "13.02"
Set-ExactText "E25" "  -4.92%  "
Set-ExactText "E26" "  +0.23%  "
Set-ExactText "D27" "2.72"
Set-ExactText "E27" "  -2.40%  "
Set-ExactText "D28" "7.25"
Set-ExactText "E28" "  -1.75%  "
Set-ExactText "E29" "  -7.61%  "
Set-ExactText "E30" "  -0.48%  "
Set-ExactText "D31" "27.07"
Set-ExactText "E31" "  +1.57%  "
Set-ExactText "E32" "  -0.01%  "
Set-ExactText "E33" "  -8.42%  "
Set-ExactText "D34" "2.54"
Set-ExactText "E34" "  -4.14%  "
Set-ExactText "E35" "  -2.63%  "
Set-ExactText "D36" "5.83"
Set-ExactText "E36" "  -2.10%  "
Set-ExactText "D37" "51.98"
Set-ExactText "E37" "  -0.82%  "
Set-ExactText "D38" "0.0₃0710"
Set-ExactText "E38" "  -4.71%  "
Set-ExactText "D39" "0.0388"
Set-ExactText "E39" "  -2.31%  "
Set-ExactText "D40" "421.88"
Set-ExactText "E40" "  -7.33%  "
Set-ExactText "D41" "0.118"
Set-ExactText "E41" "  -0.56%  "
Set-ExactText "D42" "8.29"
Set-ExactText "E42" "  -0.64%  "
Set-ExactText "B43" "Maker"
Set-ExactText "C43" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-ExactText "D43" "2.881.42"
Set-ExactText "E43" "  +0.14%  "
Set-ExactText "B44" "dogwifhat"
Set-ExactText "C44" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-ExactText "D44" "2.67"
Set-ExactText "E44" "  -11.88%  "
Set-ExactText "D45" "0.265"
Set-ExactText "E45" "  +0.42%  "
Set-ExactText "B46" "Fetch.AI"
Set-ExactText "C46" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-ExactText "D46" "2.13"
Set-ExactText "E46" "  -5.70%  "
Set-ExactText "B47" "USDe"
Set-ExactText "C47" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-ExactText "D47" "0.999"
Set-ExactText "E47" "  -0.05%  "
Set-ExactText "D48" "25.81"
Set-ExactText "E48" "  -2.78%  "
Set-ExactText "E49" "  -0.80%  "
Set-ExactText "E50" "  -7.33%  "
Set-ExactText "D51" "118.46"
